$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": the 45095.99999999999 / 24 row (row 4) is removed,
#     all following rows shift up one position and the last row's requested
#     quantity is revised from 56 to 20 ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("4:4").Delete()
$wsWeekly.Range("B8").Value2 = 20

# --- Sheet "Monthly Trend": requested quantities for the last two months change ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B3").Value2 = 16
$wsMonthly.Range("B4").Value2 = 40
